$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$failedText = "Adding and removing items from the shopping cart Failed"
$passedText = "Adding and removing items from the shopping cart Passed"

$ws.Range("A11").Value = $failedText
$ws.Range("A12").Value = $failedText
$ws.Range("A13").Value = $failedText
$ws.Range("A14").Value = $failedText
$ws.Range("A15").Value = $passedText
$ws.Range("A16").Value = $passedText
$ws.Range("A17").Value = $passedText
$ws.Range("A18").Value = $passedText
